# Apply the changes described by the ExprTestTemplate.xlsx diff:
#   1. Rename the sheet "ExprTest" -> "${newSheetName}"
#   2. Set the sheet's (Normal view) zoom to 100%
#      -> <sheetView ... zoomScaleNormal="100" .../>
#   3. Add a header/footer with left/center/right text on both the
#      header and the footer, each referencing an element of numberList.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet.
$ws.Name = '${newSheetName}'

# 2. Set the zoom level for the sheet's normal view to 100%.
$excel.ActiveWindow.Zoom = 100

# 3. Configure the odd header / odd footer, left/center/right sections.
$ps = $ws.PageSetup
$ps.LeftHeader   = 'Header Left: ${numberList[0]}'
$ps.CenterHeader = 'Header Center: ${numberList[1]}'
$ps.RightHeader  = 'Header Right: ${numberList[2]}'
$ps.LeftFooter   = 'Footer Left: ${numberList[3]}'
$ps.CenterFooter = 'Footer Center: ${numberList[4]}'
$ps.RightFooter  = 'Footer Right: ${numberList[5]}'
